$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combined")

$ws.Range("C4").Value = 0.00111633
$ws.Range("C5").Value = 0.0089401099999999994
$ws.Range("C6").Value = 0.014279099999999999
$ws.Range("C7").Value = 0.025055500000000001
$ws.Range("C8").Value = 0.039445800000000003
$ws.Range("C9").Value = 0.087393200000000004
$ws.Range("C10").Value = 0.12521599999999999
$ws.Range("C11").Value = 0.16814100000000001
$ws.Range("C12").Value = 0.22578500000000001
$ws.Range("C13").Value = 0.30410900000000002
$ws.Range("C22").Value = 0.050095099999999997
$ws.Range("D22").Value = 0.0057773299999999998
$ws.Range("E22").Value = 0.0058496399999999997
$ws.Range("C23").Value = 0.075114399999999998
$ws.Range("D23").Value = 0.017118499999999998
$ws.Range("E23").Value = 0.00066879400000000001
$ws.Range("C24").Value = 0.204151
$ws.Range("D24").Value = 0.032907699999999998
$ws.Range("E24").Value = 0.0013018699999999999
$ws.Range("C25").Value = 0.26647500000000002
$ws.Range("D25").Value = 0.075818800000000006
$ws.Range("E25").Value = 0.0081160299999999998
$ws.Range("C26").Value = 0.74996700000000005
$ws.Range("D26").Value = 0.162605
$ws.Range("E26").Value = 0.0098357800000000006
$ws.Range("C27").Value = 0.82633800000000002
$ws.Range("D27").Value = 0.21192
$ws.Range("E27").Value = 0.010606600000000001
$ws.Range("C28").Value = 1.35978
$ws.Range("D28").Value = 0.366394
$ws.Range("E28").Value = 0.010413199999999999
$ws.Range("C29").Value = 1.9461999999999999
$ws.Range("D29").Value = 0.70277599999999996
$ws.Range("E29").Value = 0.014773
$ws.Range("C30").Value = 3.3914900000000001
$ws.Range("D30").Value = 0.72885
$ws.Range("E30").Value = 0.026062200000000001
$ws.Range("C31").Value = 4.3013399999999997
$ws.Range("D31").Value = 0.88424800000000003
$ws.Range("E31").Value = 0.0215956
